# Update NATMI LR-pair stats (Psen1-Notch2) after re-running with corrected
# ligand/receptor expressing-cell counts (per Dr Hou's advice).
#
# For every data row (2-17) the "Ligand-expressing cells" (E) and
# "Receptor-expressing cells" (K) counts change from 1 to 3, and all of the
# dependent average/total expression & specificity columns (G,H,I,J,M,N,O,P)
# and the derived edge-weight/specificity columns (Q,R,S,T) are recomputed
# accordingly. F and L (detection rate columns) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{Row=2; "E"=3; "G"=27.58598166666667; "H"=82.75794500000001; "I"=0.2704460545904799; "J"=0.2704460545904799; "K"=3; "M"=27.53580066666666; "N"=82.60740199999999; "O"=0.2054887285464767; "P"=0.2054887285464768; "Q"=759.6020923676544; "R"=6836.41883130889; "S"=0.05557361589820876; "T"=0.05557361589820876},
    @{Row=3; "E"=3; "G"=27.58598166666667; "H"=82.75794500000001; "I"=0.2704460545904799; "J"=0.2704460545904799; "K"=3; "M"=27.50472933333333; "N"=82.514188; "O"=0.2052568555438283; "P"=0.2052568555438283; "Q"=758.7449591359623; "R"=6828.704632223661; "S"=0.05551090675947644; "T"=0.05551090675947643},
    @{Row=4; "E"=3; "G"=27.58598166666667; "H"=82.75794500000001; "I"=0.2704460545904799; "J"=0.2704460545904799; "K"=3; "M"=45.12975566666668; "N"=135.389267; "O"=0.3367854170582615; "P"=0.3367854170582616; "Q"=1244.948612441813; "R"=11204.53751197632; "S"=0.09108228728701616; "T"=0.09108228728701614},
    @{Row=5; "E"=3; "G"=27.58598166666667; "H"=82.75794500000001; "I"=0.2704460545904799; "J"=0.2704460545904799; "K"=3; "M"=33.831228; "N"=101.493684; "O"=0.2524689988514334; "P"=0.2524689988514334; "Q"=933.2676353688202; "R"=8399.408718319381; "S"=0.06827924464577856; "T"=0.06827924464577854},
    @{Row=6; "E"=3; "G"=26.23504533333333; "H"=78.705136; "I"=0.2572018131577233; "J"=0.2572018131577233; "K"=3; "M"=27.53580066666666; "N"=82.60740199999999; "O"=0.2054887285464767; "P"=0.2054887285464768; "Q"=722.4029787796301; "R"=6501.626809016671; "S"=0.05285207356562904; "T"=0.05285207356562904},
    @{Row=7; "E"=3; "G"=26.23504533333333; "H"=78.705136; "I"=0.2572018131577233; "J"=0.2572018131577233; "K"=3; "M"=27.50472933333333; "N"=82.514188; "O"=0.2052568555438283; "P"=0.2052568555438283; "Q"=721.5878209410631; "R"=6494.290388469568; "S"=0.05279243540892553; "T"=0.05279243540892552},
    @{Row=8; "E"=3; "G"=26.23504533333333; "H"=78.705136; "I"=0.2572018131577233; "J"=0.2572018131577233; "K"=3; "M"=45.12975566666668; "N"=135.389267; "O"=0.3367854170582615; "P"=0.3367854170582616; "Q"=1183.981185797257; "R"=10655.83067217531; "S"=0.08662181991246491; "T"=0.08662181991246491},
    @{Row=9; "E"=3; "G"=26.23504533333333; "H"=78.705136; "I"=0.2572018131577233; "J"=0.2572018131577233; "K"=3; "M"=33.831228; "N"=101.493684; "O"=0.2524689988514334; "P"=0.2524689988514334; "Q"=887.563800262336; "R"=7988.074202361024; "S"=0.06493548427070384; "T"=0.06493548427070382},
    @{Row=10; "E"=3; "G"=29.31506333333334; "H"=87.94519000000001; "I"=0.2873975381543141; "J"=0.2873975381543141; "K"=3; "M"=27.53580066666666; "N"=82.60740199999999; "O"=0.2054887285464767; "P"=0.2054887285464768; "Q"=807.2137404773756; "R"=7264.923664296381; "S"=0.05905695470271755; "T"=0.05905695470271755},
    @{Row=11; "E"=3; "G"=29.31506333333334; "H"=87.94519000000001; "I"=0.2873975381543141; "J"=0.2873975381543141; "K"=3; "M"=27.50472933333333; "N"=82.514188; "O"=0.2052568555438283; "P"=0.2052568555438283; "Q"=806.3028823728579; "R"=7256.725941355721; "S"=0.05899031497259193; "T"=0.05899031497259194},
    @{Row=12; "E"=3; "G"=29.31506333333334; "H"=87.94519000000001; "I"=0.2873975381543141; "J"=0.2873975381543141; "K"=3; "M"=45.12975566666668; "N"=135.389267; "O"=0.3367854170582615; "P"=0.3367854170582616; "Q"=1322.981645586193; "R"=11906.83481027573; "S"=0.09679129974881831; "T"=0.09679129974881834},
    @{Row=13; "E"=3; "G"=29.31506333333334; "H"=87.94519000000001; "I"=0.2873975381543141; "J"=0.2873975381543141; "K"=3; "M"=33.831228; "N"=101.493684; "O"=0.2524689988514334; "P"=0.2524689988514334; "Q"=991.7645914644403; "R"=8925.881323179961; "S"=0.07255896873018632; "T"=0.07255896873018632},
    @{Row=14; "E"=3; "G"=18.86569966666666; "H"=56.59709899999999; "I"=0.1849545940974826; "J"=0.1849545940974826; "K"=3; "M"=27.53580066666666; "N"=82.60740199999999; "O"=0.2054887285464767; "P"=0.2054887285464768; "Q"=519.482145458533; "R"=4675.339309126797; "S"=0.0380060843799214; "T"=0.0380060843799214},
    @{Row=15; "E"=3; "G"=18.86569966666666; "H"=56.59709899999999; "I"=0.1849545940974826; "J"=0.1849545940974826; "K"=3; "M"=27.50472933333333; "N"=82.514188; "O"=0.2052568555438283; "P"=0.2052568555438283; "Q"=518.8959630156235; "R"=4670.063667140612; "S"=0.03796319840283439; "T"=0.03796319840283439},
    @{Row=16; "E"=3; "G"=18.86569966666666; "H"=56.59709899999999; "I"=0.1849545940974826; "J"=0.1849545940974826; "K"=3; "M"=45.12975566666668; "N"=135.389267; "O"=0.3367854170582615; "P"=0.3367854170582616; "Q"=851.4044164373815; "R"=7662.639747936433; "S"=0.06229001010996216; "T"=0.06229001010996218},
    @{Row=17; "E"=3; "G"=18.86569966666666; "H"=56.59709899999999; "I"=0.1849545940974826; "J"=0.1849545940974826; "K"=3; "M"=33.831228; "N"=101.493684; "O"=0.2524689988514334; "P"=0.2524689988514334; "Q"=638.249786802524; "R"=5744.248081222715; "S"=0.04669530120476467; "T"=0.04669530120476467}
)

foreach ($rowEntry in $rowsData) {
    $rowNum = $rowEntry["Row"]
    foreach ($col in $rowEntry.Keys) {
        if ($col -ne "Row") {
            $ws.Range("$col$rowNum").Value = $rowEntry[$col]
        }
    }
}
